$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; existing rows 2-9 shift down to 3-10
$ws.Rows("2:2").Insert()
# The inserted row inherits the formatting of the row above (the bold header);
# reset it back to the workbook default so the new data row is unstyled, just
# like the rest of the data rows.
$ws.Range("A2:Q2").Style = "Normal"

# Populate the newly inserted row 2 with the new IPO record.
# Force the date-shaped text columns to be stored as plain text (not auto-converted
# to Excel date serials) by pre-setting the number format, then restoring the
# default "Normal" style once the text value has been written.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2024-04-30"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "제일엠앤에스"
$ws.Range("C2").Value = "KB"
$ws.Range("D2").Value = 528
$ws.Range("E2").Value = "코스닥"
$ws.Range("F2").Value = 528
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 22000
$ws.Range("N2").Value = 100

$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "2024-04-18"
$ws.Range("O2").Style = "Normal"

$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "2024-04-23"
$ws.Range("P2").Style = "Normal"

$ws.Range("Q2").Value = 1800000

# Clear any formatting the row-insert propagated onto the new row so the cells
# stay on the workbook's default (unstyled) cell format, matching the other rows.
$ws.Range("A2:Q2").Style = "Normal"

# All underwriting-market values (column C) change from "코스닥" to "KB" for the
# pre-existing rows that were shifted down to rows 3-10
$ws.Range("C3:C10").Value = "KB"
